$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 2851.375
$ws.Range("I6").Value = 3701.8333
$ws.Range("J6").Value = 300
$ws.Range("K6").Value = 11105.4999
$ws.Range("L6").Value = 900
$ws.Range("M6").Value = -10993.4999
$ws.Range("N6").Value = -1124
$ws.Range("H12").Value = 412.7143
$ws.Range("I12").Value = 407.25
$ws.Range("J12").Value = 420
$ws.Range("K12").Value = 407.25
$ws.Range("L12").Value = 420
$ws.Range("M12").Value = -237.25
$ws.Range("H17").Value = 2889.0908
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2889.0908
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 8667.2724
$ws.Range("N17").Value = -9003.2724
$ws.Range("H43").Value = 1484
$ws.Range("I43").Value = 1479
$ws.Range("J43").Value = 1499
$ws.Range("K43").Value = 1479
$ws.Range("L43").Value = 1499
$ws.Range("M43").Value = -1410
$ws.Range("N43").Value = -1637
$ws.Range("H62").Value = 4046.9285
$ws.Range("I62").Value = 3997.125
$ws.Range("J62").Value = 4113.3335
$ws.Range("K62").Value = 3997.125
$ws.Range("L62").Value = 4113.3335
$ws.Range("M62").Value = -3373.125
$ws.Range("N62").Value = -5361.3335
$ws.Range("H65").Value = 4046.9285
$ws.Range("I65").Value = 3997.125
$ws.Range("J65").Value = 4113.3335
$ws.Range("K65").Value = 19985.625
$ws.Range("L65").Value = 20566.6675
$ws.Range("M65").Value = -16865.625
$ws.Range("N65").Value = -26806.6675
$ws.Range("H132").Value = 2662.6316
$ws.Range("I132").Value = 2532.7778
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 7598.3334
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -5068.3334
$ws.Range("N132").Value = -20060
$ws.Range("H137").Value = 1952.4445
$ws.Range("I137").Value = 733
$ws.Range("J137").Value = 2562.1667
$ws.Range("K137").Value = 2199
$ws.Range("L137").Value = 7686.500100000001
$ws.Range("M137").Value = 351
$ws.Range("H138").Value = 3642.3447
$ws.Range("I138").Value = 2987
$ws.Range("J138").Value = 3665.75
$ws.Range("K138").Value = 8961
$ws.Range("L138").Value = 10997.25
$ws.Range("M138").Value = -3821
$ws.Range("N138").Value = -21277.25

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1684.0714
$ws.Range("I61").Value = 1698.6154
$ws.Range("J61").Value = 1495
$ws.Range("K61").Value = 1698.6154
$ws.Range("L61").Value = 1495
$ws.Range("M61").Value = -1486.6154
$ws.Range("N61").Value = -1919
$ws.Range("H136").Value = 1684.0714
$ws.Range("I136").Value = 1698.6154
$ws.Range("J136").Value = 1495
$ws.Range("K136").Value = 5095.8462
$ws.Range("L136").Value = 4485
$ws.Range("M136").Value = -2545.8462
$ws.Range("N136").Value = -9585

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 19150
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 19150
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 19150
$ws.Range("N29").Value = -19736
$ws.Range("H31").Value = 2206
$ws.Range("I31").Value = 1676.32
$ws.Range("J31").Value = 3151.8572
$ws.Range("K31").Value = 1676.32
$ws.Range("L31").Value = 3151.8572
$ws.Range("M31").Value = -1381.32
$ws.Range("N31").Value = -3741.8572
$ws.Range("H34").Value = 2206
$ws.Range("I34").Value = 1676.32
$ws.Range("J34").Value = 3151.8572
$ws.Range("K34").Value = 1676.32
$ws.Range("L34").Value = 3151.8572
$ws.Range("M34").Value = -1474.32
$ws.Range("N34").Value = -3555.8572
$ws.Range("H58").Value = 4020.5454
$ws.Range("I58").Value = 3745.3333
$ws.Range("J58").Value = 9800
$ws.Range("K58").Value = 3745.3333
$ws.Range("L58").Value = 9800
$ws.Range("M58").Value = -3542.3333
$ws.Range("N58").Value = -10206
$ws.Range("H132").Value = 1728.4445
$ws.Range("I132").Value = 1592.6666
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 4777.9998
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -2247.9998
$ws.Range("H134").Value = 2803.6667
$ws.Range("I134").Value = 1864.4
$ws.Range("J134").Value = 3977.75
$ws.Range("K134").Value = 5593.200000000001
$ws.Range("L134").Value = 11933.25
$ws.Range("M134").Value = -3058.200000000001
$ws.Range("H136").Value = 4020.5454
$ws.Range("I136").Value = 3745.3333
$ws.Range("J136").Value = 9800
$ws.Range("K136").Value = 11235.9999
$ws.Range("L136").Value = 29400
$ws.Range("M136").Value = -8685.999899999999
$ws.Range("N136").Value = -34500

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 389.875
$ws.Range("I7").Value = 250
$ws.Range("J7").Value = 473.8
$ws.Range("K7").Value = 750
$ws.Range("L7").Value = 1421.4
$ws.Range("M7").Value = -638
$ws.Range("H34").Value = 497.25
$ws.Range("I34").Value = 615
$ws.Range("J34").Value = 144
$ws.Range("K34").Value = 1845
$ws.Range("L34").Value = 432
$ws.Range("M34").Value = -1761
$ws.Range("N34").Value = -600
$ws.Range("H37").Value = 99950
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 99950
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 299850
$ws.Range("N37").Value = -300074
$ws.Range("H55").Value = 8972.35
$ws.Range("I55").Value = 13440.909
$ws.Range("J55").Value = 3510.7778
$ws.Range("K55").Value = 40322.727
$ws.Range("L55").Value = 10532.3334
$ws.Range("M55").Value = -40145.727
$ws.Range("N55").Value = -10886.3334
$ws.Range("H75").Value = 120
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 120
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 360
$ws.Range("N75").Value = -2356
$ws.Range("M75").ClearContents()
$ws.Range("H78").Value = 120
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 120
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 1080
$ws.Range("N78").Value = -11064
$ws.Range("M78").ClearContents()
$ws.Range("H87").Value = 300
$ws.Range("I87").Value = 300
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 900
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = 348
$ws.Range("H90").Value = 300
$ws.Range("I90").Value = 300
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 2700
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = 3540
$ws.Range("H114").Value = 210
$ws.Range("I114").Value = 210
$ws.Range("J114").Value = 0
$ws.Range("K114").Value = 630
$ws.Range("L114").Value = 0
$ws.Range("M114").Value = 2624
$ws.Range("N114").ClearContents()
$ws.Range("H117").Value = 638.6
$ws.Range("I117").Value = 769
$ws.Range("J117").Value = 606
$ws.Range("K117").Value = 2307
$ws.Range("L117").Value = 1818
$ws.Range("M117").Value = 1135
$ws.Range("N117").Value = -8702
$ws.Range("H121").Value = 898.25
$ws.Range("I121").Value = 769
$ws.Range("J121").Value = 1027.5
$ws.Range("K121").Value = 2307
$ws.Range("L121").Value = 3082.5
$ws.Range("M121").Value = -997
$ws.Range("H128").Value = 340382
$ws.Range("I128").Value = 340382
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 1021146
$ws.Range("L128").Value = 0
$ws.Range("M128").Value = -1016166

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 114.5
$ws.Range("I2").Value = 144.5
$ws.Range("J2").Value = 39.5
$ws.Range("K2").Value = 144.5
$ws.Range("L2").Value = 39.5
$ws.Range("M2").Value = -31.5
$ws.Range("N2").Value = -265.5
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H70").Value = 83336340
$ws.Range("I70").Value = 111113450
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 111113450
$ws.Range("L70").Value = 5000
$ws.Range("M70").Value = -111113180
$ws.Range("N70").Value = -5540
$ws.Range("H73").Value = 83336340
$ws.Range("I73").Value = 111113450
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 111113450
$ws.Range("L73").Value = 5000
$ws.Range("M73").Value = -111112514
$ws.Range("N73").Value = -6872
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H132").Value = 828.0909
$ws.Range("I132").Value = 842.5
$ws.Range("J132").Value = 684
$ws.Range("K132").Value = 2527.5
$ws.Range("L132").Value = 2052
$ws.Range("M132").Value = 2.5
$ws.Range("N132").Value = -7112

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 68994
$ws.Range("I7").Value = 68994
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 68994
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -68882
$ws.Range("H22").Value = 1115.6666
$ws.Range("I22").Value = 1129.1666
$ws.Range("J22").Value = 1061.6666
$ws.Range("K22").Value = 1129.1666
$ws.Range("L22").Value = 1061.6666
$ws.Range("M22").Value = -834.1666
$ws.Range("N22").Value = -1651.6666
$ws.Range("H27").Value = 1115.6666
$ws.Range("I27").Value = 1129.1666
$ws.Range("J27").Value = 1061.6666
$ws.Range("K27").Value = 1129.1666
$ws.Range("L27").Value = 1061.6666
$ws.Range("M27").Value = -1022.1666
$ws.Range("N27").Value = -1275.6666
$ws.Range("H126").Value = 68994
$ws.Range("I126").Value = 68994
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 206982
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -204512

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 6666.3335
$ws.Range("I2").Value = 5000
$ws.Range("J2").Value = 9999
$ws.Range("K2").Value = 5000
$ws.Range("L2").Value = 9999
$ws.Range("M2").Value = -4888
$ws.Range("N2").Value = -10223
$ws.Range("H54").Value = 25400
$ws.Range("I54").Value = 12000
$ws.Range("J54").Value = 28750
$ws.Range("K54").Value = 12000
$ws.Range("L54").Value = 28750
$ws.Range("M54").Value = -11480
$ws.Range("N54").Value = -29790
$ws.Range("H81").Value = 3887.6667
$ws.Range("I81").Value = 3887.6667
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 7775.3334
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -6714.3334
$ws.Range("H84").Value = 3887.6667
$ws.Range("I84").Value = 3887.6667
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 38876.667
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -33572.667
$ws.Range("H136").Value = 5130.6924
$ws.Range("I136").Value = 5112.45
$ws.Range("J136").Value = 5191.5
$ws.Range("K136").Value = 15337.35
$ws.Range("L136").Value = 15574.5
$ws.Range("M136").Value = -12787.35
$ws.Range("N136").Value = -20674.5
